$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebuild the hyperlinks so B1's hyperlink keeps pointing at
# mailto:dixit1234@ (rId2, unchanged) but now also carries a display
# text of "dixit1234@". The other hyperlinks are recreated unchanged
# so the relationship ids line back up (rId1..rId7) the same as before.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:test1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:dixit1234@", "", "", "dixit1234@")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:test2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:dixit1234@")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:dixitjani666@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:dixit1234@")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:test4@gmail.com")

# Re-applying the Hyperlink style keeps the cells on the original style
# record instead of the fresh duplicate that Hyperlinks.Add produces.
$ws.Range("A1").Style = "Hyperlink"
$ws.Range("B1").Style = "Hyperlink"
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("A4").Style = "Hyperlink"

# Add the new "dixit123" shared string and point B1 at it.
$ws.Range("B1").Value = "dixit123"

# Move the active selection from B4 to B3.
$null = $ws.Range("B3").Select()
